$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 1.3
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("R6").Value = 1.48
$ws.Range("R7").Value = 1.5
$ws.Range("R8").Value = 1.62
$ws.Range("G10").Value = 2.7
$ws.Range("I10").Value = 2.6
$ws.Range("L10").Value = 3.2
$ws.Range("AH10").Value = 8.5
$ws.Range("AZ10").Value = 41
$ws.Range("BD10").Value = 126
$ws.Range("M13").Value = 1.05
$ws.Range("O13").Value = 1.29
$ws.Range("Q13").Value = 1.9
$ws.Range("R13").Value = 1.95
$ws.Range("G16").Value = 2.62
$ws.Range("H16").Value = 3.05
$ws.Range("Z16").Value = 26
$ws.Range("AB16").Value = 18.5
$ws.Range("AC16").Value = 10.75
$ws.Range("AD16").Value = 5.4
$ws.Range("AM16").Value = 19
$ws.Range("AN16").Value = 4.75
$ws.Range("AO16").Value = 14
$ws.Range("AQ16").Value = 60
$ws.Range("AX16").Value = 13.5
$ws.Range("AY16").Value = 18.5
$ws.Range("AZ16").Value = 55
$ws.Range("G19").Value = 2.05
$ws.Range("H19").Value = 3.6
$ws.Range("I19").Value = 3.4
$ws.Range("J19").Value = 2.63
$ws.Range("AC19").Value = 13
$ws.Range("AK19").Value = 41
$ws.Range("AQ19").Value = 34
$ws.Range("AV19").Value = 41
$ws.Range("AX19").Value = 19
$ws.Range("O22").Value = 1.2
$ws.Range("P22").Value = 4.33
$ws.Range("Q22").Value = 1.67
$ws.Range("R22").Value = 2.15
$ws.Range("G24").Value = 3.1
$ws.Range("H24").Value = 3.4
$ws.Range("I24").Value = 2.15
$ws.Range("J24").Value = 3.5
$ws.Range("K24").Value = 2.38
$ws.Range("L24").Value = 2.75
$ws.Range("Q24").Value = 1.6
$ws.Range("R24").Value = 2.3
$ws.Range("U24").Value = 1.5
$ws.Range("V24").Value = 2.5
$ws.Range("X24").Value = 19
$ws.Range("Y24").Value = 12
$ws.Range("Z24").Value = 34
$ws.Range("AA24").Value = 23
$ws.Range("AB24").Value = 26
$ws.Range("AD24").Value = 7
$ws.Range("AE24").Value = 11
$ws.Range("AG24").Value = 101
$ws.Range("AH24").Value = 11
$ws.Range("AI24").Value = 13
$ws.Range("AK24").Value = 21
$ws.Range("AL24").Value = 15
$ws.Range("AN24").Value = 5.5
$ws.Range("AO24").Value = 15
$ws.Range("AP24").Value = 21
$ws.Range("AR24").Value = 51
$ws.Range("AS24").Value = 101
$ws.Range("AW24").Value = 4.5
$ws.Range("AX24").Value = 11
$ws.Range("AZ24").Value = 34
$ws.Range("BA24").Value = 51
$ws.Range("G25").Value = 2.05
$ws.Range("H25").Value = 3.25
$ws.Range("I25").Value = 3.9
$ws.Range("J25").Value = 2.75
$ws.Range("L25").Value = 4.33
$ws.Range("X25").Value = 9
$ws.Range("Y25").Value = 9
$ws.Range("Z25").Value = 17
$ws.Range("AG25").Value = 401
$ws.Range("AH25").Value = 9.5
$ws.Range("AI25").Value = 19
$ws.Range("AL25").Value = 34
$ws.Range("AY25").Value = 34
$ws.Range("AZ25").Value = 81
$ws.Range("H26").Value = 3
$ws.Range("I26").Value = 3.1
$ws.Range("L26").Value = 4
$ws.Range("M26").Value = 1.11
$ws.Range("N26").Value = 6.5
$ws.Range("AS26").Value = 301
$ws.Range("I27").Value = 2.7
$ws.Range("K27").Value = 2.05
$ws.Range("L27").Value = 3.4
$ws.Range("S27").Value = 1.44
$ws.Range("T27").Value = 2.63
$ws.Range("U27").Value = 1.91
$ws.Range("V27").Value = 1.91
$ws.Range("W27").Value = 8
$ws.Range("X27").Value = 13
$ws.Range("Z27").Value = 29
$ws.Range("AC27").Value = 8
$ws.Range("AD27").Value = 6
$ws.Range("AG27").Value = 301
$ws.Range("AK27").Value = 26
$ws.Range("AN27").Value = 4.75
$ws.Range("AP27").Value = 26
$ws.Range("AT27").Value = 2.63
$ws.Range("AU27").Value = 8
$ws.Range("AV27").Value = 51
$ws.Range("AX27").Value = 15
$ws.Range("AY27").Value = 26
$ws.Range("G31").Value = 2.15
$ws.Range("J31").Value = 2.88
$ws.Range("I32").Value = 1.55
$ws.Range("G33").Value = 1.85
$ws.Range("K33").Value = 1.91
$ws.Range("G34").Value = 1.33
$ws.Range("J34").Value = 1.95
$ws.Range("G36").Value = 1.18
$ws.Range("H36").Value = 7
$ws.Range("I36").Value = 11
$ws.Range("J36").Value = 1.57
$ws.Range("K36").Value = 2.88
$ws.Range("L36").Value = 10
$ws.Range("N36").Value = 21
$ws.Range("Q36").Value = 1.44
$ws.Range("R36").Value = 2.7
$ws.Range("W36").Value = 9
$ws.Range("Y36").Value = 10
$ws.Range("AD36").Value = 13
$ws.Range("AE36").Value = 26
$ws.Range("AF36").Value = 67
$ws.Range("AJ36").Value = 29
$ws.Range("AK36").Value = 151
$ws.Range("AQ36").Value = 12
$ws.Range("AU36").Value = 10
$ws.Range("AV36").Value = 51
$ws.Range("AW36").Value = 12
$ws.Range("AY36").Value = 41
$ws.Range("AZ36").Value = 251
$ws.Range("BA36").Value = 201
$ws.Range("G37").Value = 1.67
$ws.Range("Q37").Value = 2.5
$ws.Range("R37").Value = 1.5
$ws.Range("M39").Value = 1.11
$ws.Range("N39").Value = 6.5
$ws.Range("O41").Value = 1.33
$ws.Range("P41").Value = 3.25
$ws.Range("Q44").Value = 2.5
$ws.Range("R44").Value = 1.5
$ws.Range("G45").Value = 1.91
$ws.Range("H45").Value = 3.5
$ws.Range("I45").Value = 4
$ws.Range("J45").Value = 2.5
$ws.Range("L45").Value = 4
$ws.Range("O45").Value = 1.22
$ws.Range("P45").Value = 4
$ws.Range("Q45").Value = 1.73
$ws.Range("R45").Value = 2.08
$ws.Range("Z45").Value = 17
$ws.Range("AA45").Value = 15
$ws.Range("AB45").Value = 23
$ws.Range("AD45").Value = 6.5
$ws.Range("AE45").Value = 12
$ws.Range("AI45").Value = 21
$ws.Range("AJ45").Value = 13
$ws.Range("AL45").Value = 29
$ws.Range("AO45").Value = 10
$ws.Range("AP45").Value = 19
$ws.Range("AQ45").Value = 34
$ws.Range("AR45").Value = 51
$ws.Range("AX45").Value = 19
$ws.Range("AY45").Value = 23
$ws.Range("AZ45").Value = 51
$ws.Range("BA45").Value = 67
$ws.Range("I46").Value = 2
$ws.Range("N46").Value = 17
$ws.Range("Q46").Value = 1.6
$ws.Range("R46").Value = 2.3
$ws.Range("AB46").Value = 26
$ws.Range("AO46").Value = 17
$ws.Range("AP46").Value = 21
